$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new test case rows for the "Snake MG" mini-game (rows 20-24).
# Cell values are written in the exact order the original author entered
# them so that the shared-string table ends up in the same sequence.

$ws.Cells.Item(20, 1).Value = 2
$ws.Cells.Item(20, 2).Value = "Snake MG"
$ws.Cells.Item(20, 4).Value = "Snake changes direction from user input"
$ws.Cells.Item(20, 3).Value = "Snake Moves"
$ws.Cells.Item(20, 5).Value = "pass"

$ws.Cells.Item(21, 3).Value = "Food Pellets Appear"
$ws.Cells.Item(21, 4).Value = "Food pellets appear randomly around the play area"
$ws.Cells.Item(21, 5).Value = "pass"

$ws.Cells.Item(22, 4).Value = "Game ends if snake collides with boundaries"
$ws.Cells.Item(22, 5).Value = "pass"

$ws.Cells.Item(23, 3).Value = "Snake Collision Game Over"
$ws.Cells.Item(23, 4).Value = "Game ends if snake collides with itself"
$ws.Cells.Item(23, 5).Value = "pass"

$ws.Cells.Item(24, 3).Value = "Score Counter"
$ws.Cells.Item(24, 4).Value = "Score counter updates and accurately displays remaining pellets needed"
$ws.Cells.Item(24, 5).Value = "pass"

$ws.Cells.Item(22, 3).Value = "Out of Bounds Game Over"

# Widen column D so the longer expected-result text fits (stored width 65).
$ws.Columns.Item(4).ColumnWidth = 64.16666667

# Match the updated selection in the saved worksheet view.
$ws.Range("C23").Select() | Out-Null
